# task-list.xlsx update:
#  - add 5 new tasks to the "Tasks" sheet (inserted after the existing row 63,
#    matching where the sheet's sortState/filter block picks back up)
#  - mark the "Notify user about PanelPRO run time/wait time" task Complete
#  - re-apply the AutoFilter (blank filter on the Status column) over the
#    grown range so the hidden/visible rows recompute from the Status values
#  - append 5 fresh blank rows at the bottom of the list
#  - keep the _FilterDatabase name and the selection in sync

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# --- insert 5 new rows right after row 63 ------------------------------
for ($i = 0; $i -lt 5; $i++) {
    $ws.Rows(64).Insert()
}

# --- row 64: High / Complete -------------------------------------------
$ws.Range("A64").Value2 = "High"
$ws.Range("B64").Value2 = "Load Testing/Server Setting Optimization"
$ws.Range("C64").Value2 = "Configure server settings based on load test results"
$ws.Range("D64").Value2 = "Complete"
$ws.Rows(64).RowHeight = 30

# --- row 65: High / not complete ----------------------------------------
$ws.Range("A65").Value2 = "High"
$ws.Range("B65").Value2 = "Reports of app crashing frequently"
$ws.Range("C65").Value2 = "This like likely due to communication between the database and app server. Reducing the number of queries to the database and pushes to the database may alleviate the issue. For example, don’t get the user base more than once per session. Don't save the pedigree to the database with every tab change, etc. "
$ws.Rows(65).RowHeight = 45

# --- row 66: Mid / not complete ------------------------------------------
$ws.Range("A66").Value2 = "Mid"
$ws.Range("B66").Value2 = "PanelPRO model selection"
$ws.Range("C66").Value2 = "Don’t list CBC as a separate 18th cancer when proband is selecting pre-specified models and make it clear that if BC is included, CBC is automatically included."
$ws.Rows(66).RowHeight = 30

# --- row 67: Mid / not complete ------------------------------------------
$ws.Range("A67").Value2 = "Mid"
$ws.Range("B67").Value2 = "Bug Fix: CBC average person"
$ws.Range("C67").Value2 = "Don't show CBC risk facet plot if proband never had BC or already had CBC. If proband had BC but not CBC, ensure the average person and proband penetrances are corrct."
$ws.Rows(67).RowHeight = 30

# --- row 68: Mid / not complete ------------------------------------------
$ws.Range("A68").Value2 = "Mid"
$ws.Range("B68").Value2 = "Include genes and tumor markers on pedigreejs"
$ws.Range("C68").Value2 = "PedigreeJS has the capability to add genes and tumor markers under each node, just need to figure out how, see canrisk.org for example"
$ws.Rows(68).RowHeight = 30

# --- existing task "Notify user about PanelPRO run time/wait time" is now
#     at row 74 (was row 69) after the insert above; mark it Complete
$ws.Range("D74").Value2 = "Complete"

# --- append 5 blank rows at the bottom of the table ----------------------
$lastRow = $ws.UsedRange.Rows.Count
for ($i = 0; $i -lt 5; $i++) {
    $ws.Rows($lastRow + 1).Insert()
    $lastRow = $lastRow + 1
}

# --- re-apply the autofilter (Status column blank filter) over the new
#     range so hidden/visible rows are recomputed from column D -----------
$ws.AutoFilterMode = $false
$ws.Range("A1:E93").AutoFilter(4, @(""), 7)

# --- keep the _FilterDatabase defined name in sync with the filter range -
$wb.Names.Item("Tasks!_FilterDatabase").RefersTo = "=Tasks!`$A`$1:`$E`$93"

# --- match the recorded selection from the edit --------------------------
$ws.Range("C90").Select()
